# Travis_County_2019_bg_SVI.xlsx -- refresh the factor-analysis example
# output (re-run produced a different, but equivalent, variable ordering
# plus solver noise in the least significant digits) across the
# 'Significant Components', 'Loading Factors', 'All Refactor Variances',
# 'Final Variances' and 'Included and Excluded' sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: Significant Components ---
$ws = $wb.Worksheets.Item("Significant Components")
$ws.Cells.Item(2, 3).Value = "['QSERV' 'QEXTRCT' 'QESL' 'QHISPC' 'QEDLESHI' 'PPUNIT' 'QNOHLTH' 'PERCAP'`n 'QFHH']"
$ws.Cells.Item(3, 3).Value = "['PERCAP' 'QRICH' 'MDHSEVAL']"
$ws.Cells.Item(4, 3).Value = "['QAGEDEP' 'MEDAGE' 'QSSBEN']"
$ws.Cells.Item(6, 3).Value = "['QAGEDEP' 'QFEMLBR' 'QFEMALE']"

# --- Sheet: Loading Factors ---
$ws = $wb.Worksheets.Item("Loading Factors")
$ws.Cells.Item(2, 1).Value = "QSERV"
$ws.Cells.Item(2, 2).Value = 0.5739901508116552
$ws.Cells.Item(2, 3).Value = 0.3660124170322497
$ws.Cells.Item(2, 4).Value = -0.1657718183206
$ws.Cells.Item(2, 5).Value = 0.3037369628643269
$ws.Cells.Item(2, 6).Value = -0.0543169869376279
$ws.Cells.Item(3, 1).Value = "QEXTRCT"
$ws.Cells.Item(3, 2).Value = 0.7782826600388352
$ws.Cells.Item(3, 3).Value = 0.1356252037097591
$ws.Cells.Item(3, 4).Value = -0.02831620362618091
$ws.Cells.Item(3, 5).Value = 0.0646122658482791
$ws.Cells.Item(3, 6).Value = -0.2142235147783895
$ws.Cells.Item(4, 1).Value = "QESL"
$ws.Cells.Item(4, 2).Value = 0.7710373695882029
$ws.Cells.Item(4, 3).Value = 0.1699122223846949
$ws.Cells.Item(4, 4).Value = -0.04922328224422961
$ws.Cells.Item(4, 5).Value = 0.1763260269844859
$ws.Cells.Item(4, 6).Value = -0.2279019734554908
$ws.Cells.Item(5, 1).Value = "QHISPC"
$ws.Cells.Item(5, 2).Value = 0.8195686404388735
$ws.Cells.Item(5, 3).Value = 0.3528442961790385
$ws.Cells.Item(5, 4).Value = -0.1123103921751212
$ws.Cells.Item(5, 5).Value = 0.108351837525289
$ws.Cells.Item(5, 6).Value = -0.1312378923933365
$ws.Cells.Item(6, 1).Value = "QEDLESHI"
$ws.Cells.Item(6, 2).Value = 0.8620380822398342
$ws.Cells.Item(6, 3).Value = 0.2178973320321957
$ws.Cells.Item(6, 4).Value = 0.009137242391464871
$ws.Cells.Item(6, 5).Value = 0.189669336882269
$ws.Cells.Item(6, 6).Value = -0.1047063831520173
$ws.Cells.Item(7, 1).Value = "PPUNIT"
$ws.Cells.Item(7, 2).Value = 0.7163381389618707
$ws.Cells.Item(7, 3).Value = -0.04911535288863684
$ws.Cells.Item(7, 4).Value = -0.08826231699206404
$ws.Cells.Item(7, 5).Value = -0.3752405488729229
$ws.Cells.Item(7, 6).Value = 0.1082069702948971
$ws.Cells.Item(8, 1).Value = "QNOHLTH"
$ws.Cells.Item(8, 2).Value = 0.6646404399934442
$ws.Cells.Item(8, 3).Value = 0.4300115830629964
$ws.Cells.Item(8, 4).Value = -0.06979631896018065
$ws.Cells.Item(8, 5).Value = 0.2528269898620109
$ws.Cells.Item(8, 6).Value = -0.1342238768424573
$ws.Cells.Item(9, 2).Value = 0.5008581705796673
$ws.Cells.Item(9, 3).Value = 0.7068286996492041
$ws.Cells.Item(9, 4).Value = -0.2354656045928186
$ws.Cells.Item(9, 5).Value = 0.1849142860830265
$ws.Cells.Item(9, 6).Value = 0.08409873265453682
$ws.Cells.Item(10, 2).Value = 0.5756381626801832
$ws.Cells.Item(10, 3).Value = 0.2481999142098628
$ws.Cells.Item(10, 4).Value = -0.007302868834883252
$ws.Cells.Item(10, 5).Value = 0.08150148749273067
$ws.Cells.Item(10, 6).Value = 0.2295452275461027
$ws.Cells.Item(11, 1).Value = "QRICH"
$ws.Cells.Item(11, 2).Value = 0.2370191022580202
$ws.Cells.Item(11, 3).Value = 0.8375862129821161
$ws.Cells.Item(11, 4).Value = -0.2059104187482207
$ws.Cells.Item(11, 5).Value = 0.3148028642276947
$ws.Cells.Item(11, 6).Value = -0.03073101708020622
$ws.Cells.Item(12, 1).Value = "MDHSEVAL"
$ws.Cells.Item(12, 2).Value = 0.3768770870459165
$ws.Cells.Item(12, 3).Value = 0.7947217985220195
$ws.Cells.Item(12, 4).Value = -0.07404588281102142
$ws.Cells.Item(12, 5).Value = -0.03403758640670176
$ws.Cells.Item(12, 6).Value = 0.02917640986609455
$ws.Cells.Item(13, 1).Value = "QAGEDEP"
$ws.Cells.Item(13, 2).Value = -0.01911179298704222
$ws.Cells.Item(13, 3).Value = -0.1421335473962037
$ws.Cells.Item(13, 4).Value = 0.7228043116834884
$ws.Cells.Item(13, 5).Value = -0.0689393639759782
$ws.Cells.Item(13, 6).Value = 0.5773032518146842
$ws.Cells.Item(14, 1).Value = "MEDAGE"
$ws.Cells.Item(14, 2).Value = -0.2799970097319023
$ws.Cells.Item(14, 3).Value = -0.2283516415427634
$ws.Cells.Item(14, 4).Value = 0.7796766711354742
$ws.Cells.Item(14, 5).Value = -0.3080563893339741
$ws.Cells.Item(14, 6).Value = -0.07828734184485212
$ws.Cells.Item(15, 2).Value = 0.03069875376985247
$ws.Cells.Item(15, 3).Value = -0.06832501896212509
$ws.Cells.Item(15, 4).Value = 0.8135614545156136
$ws.Cells.Item(15, 5).Value = -0.1401439019967177
$ws.Cells.Item(15, 6).Value = 0.1110613093848064
$ws.Cells.Item(16, 2).Value = -0.02252423356151127
$ws.Cells.Item(16, 3).Value = 0.2520805303380479
$ws.Cells.Item(16, 4).Value = -0.4372443913909764
$ws.Cells.Item(16, 5).Value = 0.737637032484032
$ws.Cells.Item(16, 6).Value = -0.1126483030971852
$ws.Cells.Item(17, 2).Value = 0.1338899665208283
$ws.Cells.Item(17, 3).Value = 0.05329239936709079
$ws.Cells.Item(17, 4).Value = -0.04656029485644806
$ws.Cells.Item(17, 5).Value = 0.7398553518884753
$ws.Cells.Item(17, 6).Value = 0.01303061387221414
$ws.Cells.Item(18, 2).Value = 0.2772247077757998
$ws.Cells.Item(18, 3).Value = 0.1320875700271804
$ws.Cells.Item(18, 4).Value = -0.3307933605032718
$ws.Cells.Item(18, 5).Value = 0.5596088513302778
$ws.Cells.Item(18, 6).Value = 0.1365168941142421
$ws.Cells.Item(19, 1).Value = "QFEMLBR"
$ws.Cells.Item(19, 2).Value = -0.2223644605289088
$ws.Cells.Item(19, 3).Value = 0.08523116655993161
$ws.Cells.Item(19, 4).Value = -0.04130466453781202
$ws.Cells.Item(19, 5).Value = 0.00490438324707819
$ws.Cells.Item(19, 6).Value = 0.7515631009125495
$ws.Cells.Item(20, 1).Value = "QFEMALE"
$ws.Cells.Item(20, 2).Value = -0.02635136301293945
$ws.Cells.Item(20, 3).Value = -0.04008883949206845
$ws.Cells.Item(20, 4).Value = 0.2036171322932024
$ws.Cells.Item(20, 5).Value = 0.02333698040035049
$ws.Cells.Item(20, 6).Value = 0.8554436033195556

# --- Sheet: All Refactor Variances ---
$ws = $wb.Worksheets.Item("All Refactor Variances")
$ws.Cells.Item(2, 9).Value = 4.83110872604232
$ws.Cells.Item(2, 10).Value = 2.779550285044058
$ws.Cells.Item(2, 11).Value = 2.319422888926871
$ws.Cells.Item(2, 12).Value = 2.298985923153931
$ws.Cells.Item(2, 13).Value = 1.901588935700678
$ws.Cells.Item(2, 14).Value = 4.904943217867458
$ws.Cells.Item(2, 15).Value = 2.604591474288087
$ws.Cells.Item(2, 16).Value = 2.297607989795691
$ws.Cells.Item(2, 17).Value = 2.045575231653374
$ws.Cells.Item(2, 18).Value = 1.900121953565291
$ws.Cells.Item(3, 9).Value = 0.241555436302116
$ws.Cells.Item(3, 10).Value = 0.1389775142522029
$ws.Cells.Item(3, 11).Value = 0.1159711444463435
$ws.Cells.Item(3, 12).Value = 0.1149492961576966
$ws.Cells.Item(3, 13).Value = 0.0950794467850339
$ws.Cells.Item(3, 14).Value = 0.2581549062035504
$ws.Cells.Item(3, 15).Value = 0.1370837618046362
$ws.Cells.Item(3, 16).Value = 0.1209267363050364
$ws.Cells.Item(3, 17).Value = 0.107661854297546
$ws.Cells.Item(3, 18).Value = 0.1000064186086995
$ws.Cells.Item(4, 9).Value = 0.241555436302116
$ws.Cells.Item(4, 10).Value = 0.3805329505543189
$ws.Cells.Item(4, 11).Value = 0.4965040950006624
$ws.Cells.Item(4, 12).Value = 0.6114533911583591
$ws.Cells.Item(4, 13).Value = 0.7065328379433929
$ws.Cells.Item(4, 14).Value = 0.2581549062035504
$ws.Cells.Item(4, 15).Value = 0.3952386680081866
$ws.Cells.Item(4, 16).Value = 0.516165404313223
$ws.Cells.Item(4, 17).Value = 0.623827258610769
$ws.Cells.Item(4, 18).Value = 0.7238336772194686
$ws.Cells.Item(5, 9).Value = 0.341888477547974
$ws.Cells.Item(5, 10).Value = 0.1967035455234392
$ws.Cells.Item(5, 11).Value = 0.1641411951692401
$ws.Cells.Item(5, 12).Value = 0.1626949095420618
$ws.Cells.Item(5, 13).Value = 0.1345718722172848
$ws.Cells.Item(5, 14).Value = 0.3566494822335782
$ws.Cells.Item(5, 15).Value = 0.1893857195636836
$ws.Cells.Item(5, 16).Value = 0.1670642581449979
$ws.Cells.Item(5, 17).Value = 0.1487383879555284
$ws.Cells.Item(5, 18).Value = 0.1381621521022118

# --- Sheet: Final Variances ---
$ws = $wb.Worksheets.Item("Final Variances")
$ws.Cells.Item(2, 2).Value = 4.904943217867458
$ws.Cells.Item(2, 3).Value = 2.604591474288087
$ws.Cells.Item(2, 4).Value = 2.297607989795691
$ws.Cells.Item(2, 5).Value = 2.045575231653374
$ws.Cells.Item(2, 6).Value = 1.900121953565291
$ws.Cells.Item(3, 2).Value = 0.2581549062035504
$ws.Cells.Item(3, 3).Value = 0.1370837618046362
$ws.Cells.Item(3, 4).Value = 0.1209267363050364
$ws.Cells.Item(3, 5).Value = 0.107661854297546
$ws.Cells.Item(3, 6).Value = 0.1000064186086995
$ws.Cells.Item(4, 2).Value = 0.2581549062035504
$ws.Cells.Item(4, 3).Value = 0.3952386680081866
$ws.Cells.Item(4, 4).Value = 0.516165404313223
$ws.Cells.Item(4, 5).Value = 0.623827258610769
$ws.Cells.Item(4, 6).Value = 0.7238336772194686
$ws.Cells.Item(5, 2).Value = 0.3566494822335782
$ws.Cells.Item(5, 3).Value = 0.1893857195636836
$ws.Cells.Item(5, 4).Value = 0.1670642581449979
$ws.Cells.Item(5, 5).Value = 0.1487383879555284
$ws.Cells.Item(5, 6).Value = 0.1381621521022118

# --- Sheet: Included and Excluded ---
$ws = $wb.Worksheets.Item("Included and Excluded")
$ws.Cells.Item(2, 2).Value = "[['QSERV', 'QEXTRCT', 'QESL', 'QHISPC', 'QEDLESHI', 'PPUNIT', 'QNOHLTH', 'PERCAP', 'QFHH', 'QRICH', 'MDHSEVAL', 'QAGEDEP', 'MEDAGE', 'QSSBEN', 'QRENTER', 'QNOAUTO', 'QPOVTY', 'QFEMLBR', 'QFEMALE']]"
